$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-25 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-26 Friday", 2)

$d.Content.Find.Execute("348×3=", $true, $false, $false, $false, $false, $true, 1, $false, "910×8=", 2)
$d.Content.Find.Execute("792×6=", $true, $false, $false, $false, $false, $true, 1, $false, "390×5=", 2)
$d.Content.Find.Execute("133×9=", $true, $false, $false, $false, $false, $true, 1, $false, "338×3=", 2)
$d.Content.Find.Execute("332×7=", $true, $false, $false, $false, $false, $true, 1, $false, "349×8=", 2)
$d.Content.Find.Execute("669×2=", $true, $false, $false, $false, $false, $true, 1, $false, "826×4=", 2)

$d.Content.Find.Execute("248×3=", $true, $false, $false, $false, $false, $true, 1, $false, "477×3=", 2)
$d.Content.Find.Execute("184×4=", $true, $false, $false, $false, $false, $true, 1, $false, "109×8=", 2)
$d.Content.Find.Execute("785×4=", $true, $false, $false, $false, $false, $true, 1, $false, "120×9=", 2)
$d.Content.Find.Execute("346×2=", $true, $false, $false, $false, $false, $true, 1, $false, "488×5=", 2)
$d.Content.Find.Execute("649×3=", $true, $false, $false, $false, $false, $true, 1, $false, "176×8=", 2)

$d.Content.Find.Execute("205×8=", $true, $false, $false, $false, $false, $true, 1, $false, "208×4=", 2)
$d.Content.Find.Execute("590×4=", $true, $false, $false, $false, $false, $true, 1, $false, "241×7=", 2)
$d.Content.Find.Execute("792×9=", $true, $false, $false, $false, $false, $true, 1, $false, "194×9=", 2)
$d.Content.Find.Execute("887×4=", $true, $false, $false, $false, $false, $true, 1, $false, "194×2=", 2)
$d.Content.Find.Execute("599×5=", $true, $false, $false, $false, $false, $true, 1, $false, "805×2=", 2)

$d.Content.Find.Execute("409×3=", $true, $false, $false, $false, $false, $true, 1, $false, "825×8=", 2)
$d.Content.Find.Execute("601×9=", $true, $false, $false, $false, $false, $true, 1, $false, "163×9=", 2)
$d.Content.Find.Execute("307×4=", $true, $false, $false, $false, $false, $true, 1, $false, "570×5=", 2)
$d.Content.Find.Execute("548×2=", $true, $false, $false, $false, $false, $true, 1, $false, "334×9=", 2)
$d.Content.Find.Execute("388×3=", $true, $false, $false, $false, $false, $true, 1, $false, "419×2=", 2)

$d.Content.Find.Execute("539×8=", $true, $false, $false, $false, $false, $true, 1, $false, "497×2=", 2)
$d.Content.Find.Execute("766×5=", $true, $false, $false, $false, $false, $true, 1, $false, "478×5=", 2)
$d.Content.Find.Execute("451×5=", $true, $false, $false, $false, $false, $true, 1, $false, "386×8=", 2)
$d.Content.Find.Execute("193×4=", $true, $false, $false, $false, $false, $true, 1, $false, "686×5=", 2)
$d.Content.Find.Execute("422×7=", $true, $false, $false, $false, $false, $true, 1, $false, "965×3=", 2)
